# Populate the strikezone/pitch-by-pitch visual data for the hitter report.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- At-bat block starting row 10 (Inning 2) ---
$ws.Range("F10").Value = "CB"
$ws.Range("G10").Value = "Take"
$ws.Range("H10").Value = "Strike"

$ws.Range("F11").Value = "CB"
$ws.Range("G11").Value = "Take"
$ws.Range("H11").Value = "Strike"

$ws.Range("F12").Value = "FB"
$ws.Range("G12").Value = "Swing"
$ws.Range("H12").Value = "Strike"
$ws.Range("M12").Value = ""

$ws.Range("M15").Value = "Strikeout"

$ws.Range("J17").Value = "CH,CB,FB"

# --- At-bat block starting row 19 (Inning 3) ---
$ws.Range("F19").Value = "CB"
$ws.Range("G19").Value = "Take"
$ws.Range("H19").Value = "Ball"
$ws.Range("M19").Value = "89.23 MPH"

$ws.Range("F20").Value = "FB"
$ws.Range("G20").Value = "Swing"
$ws.Range("H20").Value = "In Play"

$ws.Range("M21").Value = "22.91°"

$ws.Range("J26").Value = "CH,CB,FB"

# --- At-bat block starting row 28 ---
$ws.Range("F28").Value = "CH"
$ws.Range("G28").Value = "Swing"
$ws.Range("H28").Value = "In Play"
$ws.Range("M28").Value = "88.36 MPH"

$ws.Range("M30").Value = "49.18°"

$ws.Range("J35").Value = "CH,CB,FB,SL"

# --- At-bat block starting row 37 ---
$ws.Range("F37").Value = "CH"
$ws.Range("G37").Value = "Take"
$ws.Range("H37").Value = "Ball"

$ws.Range("F38").Value = "CH"
$ws.Range("G38").Value = "Swing"
$ws.Range("H38").Value = "Strike"

$ws.Range("F39").Value = "CH"
$ws.Range("G39").Value = "Swing"
$ws.Range("H39").Value = "Foul"
$ws.Range("M39").Value = ""

$ws.Range("F40").Value = "SL"
$ws.Range("G40").Value = "Swing"
$ws.Range("H40").Value = "Strike"

$ws.Range("M42").Value = "Strikeout"

$ws.Range("J44").Value = "CH,CB,FB,SL"

# --- At-bat block starting row 46 ---
$ws.Range("F46").Value = "FB"
$ws.Range("G46").Value = "Take"
$ws.Range("H46").Value = "HBP"

$ws.Range("M48").Value = ""

$ws.Range("M51").Value = "Hit By Pitch"

$ws.Range("J53").Value = "CH,FB,SL"

# --- At-bat block starting row 61 ---
$ws.Range("F61").Value = "CH"
$ws.Range("G61").Value = "Swing"
$ws.Range("H61").Value = "Foul"

$ws.Range("F62").Value = "CH"
$ws.Range("G62").Value = "Take"
$ws.Range("H62").Value = "Ball"

$ws.Range("F63").Value = "CH"
$ws.Range("G63").Value = "Swing"
$ws.Range("H63").Value = "Foul"
$ws.Range("M63").Value = ""

$ws.Range("F64").Value = "CH"
$ws.Range("G64").Value = "Take"
$ws.Range("H64").Value = "Ball"

$ws.Range("F65").Value = "CH"
$ws.Range("G65").Value = "Swing"
$ws.Range("H65").Value = "Foul"

$ws.Range("F66").Value = "CH"
$ws.Range("G66").Value = "Take"
$ws.Range("H66").Value = "Ball"
$ws.Range("M66").Value = "Strikeout"

$ws.Range("F67").Value = "CH"
$ws.Range("G67").Value = "Swing"
$ws.Range("H67").Value = "Foul"

$ws.Range("F68").Value = "CH"
$ws.Range("G68").Value = "Swing"
$ws.Range("H68").Value = "Strike"
$ws.Range("J68").Value = "CH,FB,SL"
